$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "unit/year" -> "seconds" everywhere it is used as a shared string (column G
# label on the duration/event rows). Update the cells directly so the shared
# string table is rewritten.
$ws.Range("G7").Value  = "seconds"
$ws.Range("G9").Value  = "seconds"
$ws.Range("G11").Value = "seconds"
$ws.Range("G13").Value = "seconds"

# Fix the duration formula: 340 days -> 30 days.
$ws.Range("F7").Formula  = "=30*24*60*60"
$ws.Range("F9").Formula  = "=30*24*60*60"
$ws.Range("F11").Formula = "=30*24*60*60"
$ws.Range("F13").Formula = "=30*24*60*60"

# Fix the little bug: rates were off by a factor of 10.
$ws.Range("F8").Value  = 0.6
$ws.Range("F10").Value = 0.5
$ws.Range("F12").Value = 0.1

# Move the current selection from F13 to D7 (and drop the scrolled
# topLeftCell pin that came with it).
$ws.Range("D7").Select()
